# Add Jacob Felten's Thursday hours (2) to the "Week 1" timesheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 1")
$ws.Activate()
$ws.Range("B6").Value = 2
$ws.Range("B6").Select()
